$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.845.33"
Set-TextValue $ws.Range("E2") "  +1.56%  "
Set-TextValue $ws.Range("D3") "1.763.80"
Set-TextValue $ws.Range("E3") "  +1.73%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "327.89"
Set-TextValue $ws.Range("E5") "  +1.73%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.4449"
Set-TextValue $ws.Range("E7") "  -1.69%  "
Set-TextValue $ws.Range("D8") "0.3543"
Set-TextValue $ws.Range("E8") "  +0.48%  "
Set-TextValue $ws.Range("D9") "41.96"
Set-TextValue $ws.Range("E9") "  +1.64%  "
Set-TextValue $ws.Range("D10") "0.07383"
Set-TextValue $ws.Range("E10") "  +0.42%  "
Set-TextValue $ws.Range("D11") "1.097"
Set-TextValue $ws.Range("E11") "  +2.31%  "
Set-TextValue $ws.Range("D12") "1.002"
Set-TextValue $ws.Range("E12") "  +0.04%  "
Set-TextValue $ws.Range("D13") "20.89"
Set-TextValue $ws.Range("E13") "  +2.59%  "
Set-TextValue $ws.Range("D14") "6.011"
Set-TextValue $ws.Range("E14") "  +1.82%  "
Set-TextValue $ws.Range("D15") "7.222"
Set-TextValue $ws.Range("E15") "  +2.24%  "
Set-TextValue $ws.Range("D16") "1.761.56"
Set-TextValue $ws.Range("E16") "  +2.02%  "
Set-TextValue $ws.Range("D17") "93.08"
Set-TextValue $ws.Range("E17") "  +2.52%  "
Set-TextValue $ws.Range("D18") "0.00001059"
Set-TextValue $ws.Range("E18") "  +0.85%  "
Set-TextValue $ws.Range("D19") "0.06420"
Set-TextValue $ws.Range("E19") "  +1.43%  "
Set-TextValue $ws.Range("E20") "  +0.13%  "
Set-TextValue $ws.Range("D21") "17.10"
Set-TextValue $ws.Range("E21") "  +3.18%  "
Set-TextValue $ws.Range("D22") "5.760"
Set-TextValue $ws.Range("E22") "  +0.41%  "
Set-TextValue $ws.Range("D23") "27.885.53"
Set-TextValue $ws.Range("E23") "  +1.51%  "
Set-TextValue $ws.Range("D24") "11.24"
Set-TextValue $ws.Range("E24") "  +1.09%  "
Set-TextValue $ws.Range("D25") "2.107"
Set-TextValue $ws.Range("E25") "  +1.60%  "
Set-TextValue $ws.Range("D26") "161.86"
Set-TextValue $ws.Range("E26") "  +0.36%  "
Set-TextValue $ws.Range("D27") "20.34"
Set-TextValue $ws.Range("E27") "  +2.30%  "
Set-TextValue $ws.Range("D28") "1.965.97"
Set-TextValue $ws.Range("E28") "  +1.96%  "
Set-TextValue $ws.Range("D29") "2.148"
Set-TextValue $ws.Range("E29") "  +5.15%  "
Set-TextValue $ws.Range("D30") "125.06"
Set-TextValue $ws.Range("E30") "  +0.29%  "
Set-TextValue $ws.Range("D31") "1.103"
Set-TextValue $ws.Range("E31") "  +6.09%  "
Set-TextValue $ws.Range("D32") "0.09188"
Set-TextValue $ws.Range("E32") "  +0.85%  "
Set-TextValue $ws.Range("D33") "3.685"
Set-TextValue $ws.Range("E33") "  +0.90%  "
Set-TextValue $ws.Range("D34") "5.611"
Set-TextValue $ws.Range("E34") "  +4.30%  "
Set-TextValue $ws.Range("D35") "11.82"
Set-TextValue $ws.Range("E35") "  +2.35%  "
Set-TextValue $ws.Range("D36") "0.06172"
Set-TextValue $ws.Range("E36") "  +4.00%  "
Set-TextValue $ws.Range("D37") "0.02281"
Set-TextValue $ws.Range("E37") "  +0.91%  "
Set-TextValue $ws.Range("D38") "0.2096"
Set-TextValue $ws.Range("E38") "  +2.72%  "
Set-TextValue $ws.Range("D39") "0.6302"
Set-TextValue $ws.Range("E39") "  +1.37%  "
Set-TextValue $ws.Range("D40") "4.952"
Set-TextValue $ws.Range("E40") "  +2.08%  "
Set-TextValue $ws.Range("D41") "1.187"
Set-TextValue $ws.Range("E41") "  -0.26%  "
Set-TextValue $ws.Range("D42") "1.392"
Set-TextValue $ws.Range("E42") "  +1.54%  "
Set-TextValue $ws.Range("D43") "7.868"
Set-TextValue $ws.Range("E43") "  +2.33%  "
Set-TextValue $ws.Range("D44") "13.16"
Set-TextValue $ws.Range("E44") "  +0.88%  "
Set-TextValue $ws.Range("D45") "3.749"
Set-TextValue $ws.Range("E45") "  +1.45%  "
Set-TextValue $ws.Range("D46") "0.5856"
Set-TextValue $ws.Range("E46") "  +1.41%  "
Set-TextValue $ws.Range("D47") "122.32"
Set-TextValue $ws.Range("E47") "  +0.53%  "
Set-TextValue $ws.Range("D48") "1.951"
Set-TextValue $ws.Range("E48") "  +1.70%  "
Set-TextValue $ws.Range("D49") "0.06897"
Set-TextValue $ws.Range("E49") "  +0.92%  "
Set-TextValue $ws.Range("D50") "1.132"
Set-TextValue $ws.Range("E50") "  +2.11%  "
Set-TextValue $ws.Range("D51") "72.75"
Set-TextValue $ws.Range("E51") "  +2.87%  "
